$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column-A date-like strings to be written as literal text (matching the
# original workbook, where these cells are shared-string text, not date serials).
# Setting NumberFormat "@" before the assignment prevents Excel's automatic
# "looks like a date" parsing; resetting the Style back to "Normal" afterwards
# drops the temporary number-format styling so the cell ends up unstyled, just
# like the surrounding rows.
$dateRange = $ws.Range("A147:A171")
$dateRange.NumberFormat = "@"

$ws.Range("A147").Value = "02-08-2021"
$ws.Range("B147").Value = 327.4
$ws.Range("C147").Value = 397
$ws.Range("D147").Value = 211
$ws.Range("E147").Value = 322
$ws.Range("F147").Value = 94.40000000000001
$ws.Range("G147").Value = 126.2
$ws.Range("H147").Value = 20.9
$ws.Range("I147").Value = 176.8
$ws.Range("J147").Value = 477.5
$ws.Range("K147").Value = 1594.9
$ws.Range("L147").Value = 283
$ws.Range("M147").Value = 276
$ws.Range("N147").Value = 146
$ws.Range("O147").Value = 366
$ws.Range("P147").Value = 193

$ws.Range("A148").Value = "03-08-2021"
$ws.Range("B148").Value = 326.4
$ws.Range("C148").Value = 396
$ws.Range("D148").Value = 209
$ws.Range("E148").Value = 321
$ws.Range("F148").Value = 93.8
$ws.Range("G148").Value = 125.9
$ws.Range("H148").Value = 20.8
$ws.Range("I148").Value = 174.2
$ws.Range("J148").Value = 477.9
$ws.Range("K148").Value = 1594.7
$ws.Range("L148").Value = 284.3
$ws.Range("M148").Value = 276
$ws.Range("N148").Value = 144
$ws.Range("O148").Value = 367
$ws.Range("P148").Value = 187

$ws.Range("A149").Value = "04-08-2021"
$ws.Range("B149").Value = 325.8
$ws.Range("C149").Value = 396
$ws.Range("D149").Value = 208
$ws.Range("E149").Value = 321
$ws.Range("F149").Value = 93.8
$ws.Range("G149").Value = 127.3
$ws.Range("H149").Value = 19.8
$ws.Range("I149").Value = 170.9
$ws.Range("J149").Value = 481.8
$ws.Range("K149").Value = 1589.1
$ws.Range("L149").Value = 284.7
$ws.Range("M149").Value = 277
$ws.Range("N149").Value = 144
$ws.Range("O149").Value = 367
$ws.Range("P149").Value = 186

$ws.Range("A150").Value = "05-08-2021"
$ws.Range("B150").Value = 322.9
$ws.Range("C150").Value = 393
$ws.Range("D150").Value = 206
$ws.Range("E150").Value = 319
$ws.Range("F150").Value = 92.90000000000001
$ws.Range("G150").Value = 126.1
$ws.Range("H150").Value = 18.7
$ws.Range("I150").Value = 165.8
$ws.Range("J150").Value = 483.3
$ws.Range("K150").Value = 1580.2
$ws.Range("L150").Value = 280.9
$ws.Range("M150").Value = 276
$ws.Range("N150").Value = 140
$ws.Range("O150").Value = 364
$ws.Range("P150").Value = 182

$ws.Range("A151").Value = "06-08-2021"
$ws.Range("B151").Value = 319.4
$ws.Range("C151").Value = 389
$ws.Range("D151").Value = 203
$ws.Range("E151").Value = 315
$ws.Range("F151").Value = 92.59999999999999
$ws.Range("G151").Value = 124.6
$ws.Range("H151").Value = 16.9
$ws.Range("I151").Value = 160.4
$ws.Range("J151").Value = 482.4
$ws.Range("K151").Value = 1572.8
$ws.Range("L151").Value = 280.2
$ws.Range("M151").Value = 274
$ws.Range("N151").Value = 136
$ws.Range("O151").Value = 360
$ws.Range("P151").Value = 179

$ws.Range("A152").Value = "09-08-2021"
$ws.Range("B152").Value = 319.3
$ws.Range("C152").Value = 391
$ws.Range("D152").Value = 203
$ws.Range("E152").Value = 315
$ws.Range("F152").Value = 92.90000000000001
$ws.Range("G152").Value = 123.7
$ws.Range("H152").Value = 15.7
$ws.Range("I152").Value = 157.5
$ws.Range("J152").Value = 485
$ws.Range("K152").Value = 1578.5
$ws.Range("L152").Value = 280
$ws.Range("M152").Value = 275
$ws.Range("N152").Value = 135
$ws.Range("O152").Value = 363
$ws.Range("P152").Value = 179

$ws.Range("A153").Value = "10-08-2021"
$ws.Range("B153").Value = 318.5
$ws.Range("C153").Value = 391
$ws.Range("D153").Value = 202
$ws.Range("E153").Value = 312
$ws.Range("F153").Value = 92.8
$ws.Range("G153").Value = 123.6
$ws.Range("H153").Value = 13.9
$ws.Range("I153").Value = 154.9
$ws.Range("J153").Value = 482.2
$ws.Range("K153").Value = 1564.3
$ws.Range("L153").Value = 279.9
$ws.Range("M153").Value = 277
$ws.Range("N153").Value = 136
$ws.Range("O153").Value = 362
$ws.Range("P153").Value = 183

$ws.Range("A154").Value = "11-08-2021"
$ws.Range("B154").Value = 319.8
$ws.Range("C154").Value = 391
$ws.Range("D154").Value = 204
$ws.Range("E154").Value = 315
$ws.Range("F154").Value = 92.09999999999999
$ws.Range("G154").Value = 120.4
$ws.Range("H154").Value = 17.3
$ws.Range("I154").Value = 157.1
$ws.Range("J154").Value = 486.1
$ws.Range("K154").Value = 1564.5
$ws.Range("L154").Value = 280.3
$ws.Range("M154").Value = 279
$ws.Range("N154").Value = 136
$ws.Range("O154").Value = 359
$ws.Range("P154").Value = 184

$ws.Range("A155").Value = "12-08-2021"
$ws.Range("B155").Value = 317.1
$ws.Range("C155").Value = 387
$ws.Range("D155").Value = 201
$ws.Range("E155").Value = 311
$ws.Range("F155").Value = 91.7
$ws.Range("G155").Value = 118.7
$ws.Range("H155").Value = 12.9
$ws.Range("I155").Value = 154.6
$ws.Range("J155").Value = 478.2
$ws.Range("K155").Value = 1538.3
$ws.Range("L155").Value = 279.1
$ws.Range("M155").Value = 277
$ws.Range("N155").Value = 135
$ws.Range("O155").Value = 355
$ws.Range("P155").Value = 183

$ws.Range("A156").Value = "13-08-2021"
$ws.Range("B156").Value = 321.7
$ws.Range("C156").Value = 392
$ws.Range("D156").Value = 205
$ws.Range("E156").Value = 313
$ws.Range("F156").Value = 91.5
$ws.Range("G156").Value = 118.5
$ws.Range("H156").Value = 13.4
$ws.Range("I156").Value = 160.8
$ws.Range("J156").Value = 475.9
$ws.Range("K156").Value = 1566.1
$ws.Range("L156").Value = 282.7
$ws.Range("M156").Value = 279
$ws.Range("N156").Value = 139
$ws.Range("O156").Value = 359
$ws.Range("P156").Value = 188

$ws.Range("A157").Value = "16-08-2021"
$ws.Range("B157").Value = 323.5
$ws.Range("C157").Value = 396
$ws.Range("D157").Value = 205
$ws.Range("E157").Value = 314
$ws.Range("F157").Value = 92.09999999999999
$ws.Range("G157").Value = 118.9
$ws.Range("H157").Value = 15.9
$ws.Range("I157").Value = 162
$ws.Range("J157").Value = 474.6
$ws.Range("K157").Value = 1589.1
$ws.Range("L157").Value = 285.9
$ws.Range("M157").Value = 281
$ws.Range("N157").Value = 141
$ws.Range("O157").Value = 364
$ws.Range("P157").Value = 189

$ws.Range("A158").Value = "17-08-2021"
$ws.Range("B158").Value = 322.9
$ws.Range("C158").Value = 396
$ws.Range("D158").Value = 205
$ws.Range("E158").Value = 312
$ws.Range("F158").Value = 91.59999999999999
$ws.Range("G158").Value = 120
$ws.Range("H158").Value = 14.1
$ws.Range("I158").Value = 159.7
$ws.Range("J158").Value = 473.3
$ws.Range("K158").Value = 1589
$ws.Range("L158").Value = 287
$ws.Range("M158").Value = 280
$ws.Range("N158").Value = 141
$ws.Range("O158").Value = 364
$ws.Range("P158").Value = 188

$ws.Range("A159").Value = "18-08-2021"
$ws.Range("B159").Value = 323
$ws.Range("C159").Value = 395
$ws.Range("D159").Value = 204
$ws.Range("E159").Value = 313
$ws.Range("F159").Value = 91.2
$ws.Range("G159").Value = 120
$ws.Range("H159").Value = 15
$ws.Range("I159").Value = 159.7
$ws.Range("J159").Value = 475.8
$ws.Range("K159").Value = 1589.4
$ws.Range("L159").Value = 290.7
$ws.Range("M159").Value = 280
$ws.Range("N159").Value = 140
$ws.Range("O159").Value = 364
$ws.Range("P159").Value = 187

$ws.Range("A160").Value = "19-08-2021"
$ws.Range("B160").Value = 325.4
$ws.Range("C160").Value = 399
$ws.Range("D160").Value = 204
$ws.Range("E160").Value = 317
$ws.Range("F160").Value = 91.3
$ws.Range("G160").Value = 119.9
$ws.Range("H160").Value = 15.6
$ws.Range("I160").Value = 161.9
$ws.Range("J160").Value = 481.1
$ws.Range("K160").Value = 1603.6
$ws.Range("L160").Value = 297.9
$ws.Range("M160").Value = 285
$ws.Range("N160").Value = 142
$ws.Range("O160").Value = 369
$ws.Range("P160").Value = 187

$ws.Range("A161").Value = "20-08-2021"
$ws.Range("B161").Value = 324
$ws.Range("C161").Value = 398
$ws.Range("D161").Value = 203
$ws.Range("E161").Value = 315
$ws.Range("F161").Value = 90.3
$ws.Range("G161").Value = 119.4
$ws.Range("H161").Value = 14.9
$ws.Range("I161").Value = 160.7
$ws.Range("J161").Value = 478.4
$ws.Range("K161").Value = 1597.9
$ws.Range("L161").Value = 297.3
$ws.Range("M161").Value = 286
$ws.Range("N161").Value = 141
$ws.Range("O161").Value = 369
$ws.Range("P161").Value = 185

$ws.Range("A162").Value = "23-08-2021"
$ws.Range("B162").Value = 323.7
$ws.Range("C162").Value = 396
$ws.Range("D162").Value = 204
$ws.Range("E162").Value = 316
$ws.Range("F162").Value = 91.2
$ws.Range("G162").Value = 119.8
$ws.Range("H162").Value = 16.3
$ws.Range("I162").Value = 162.2
$ws.Range("J162").Value = 479.7
$ws.Range("K162").Value = 1571.9
$ws.Range("L162").Value = 295.5
$ws.Range("M162").Value = 286
$ws.Range("N162").Value = 139
$ws.Range("O162").Value = 368
$ws.Range("P162").Value = 185

$ws.Range("A163").Value = "24-08-2021"
$ws.Range("B163").Value = 318.9
$ws.Range("C163").Value = 390
$ws.Range("D163").Value = 202
$ws.Range("E163").Value = 309
$ws.Range("F163").Value = 90.40000000000001
$ws.Range("G163").Value = 120.1
$ws.Range("H163").Value = 11.5
$ws.Range("I163").Value = 158.2
$ws.Range("J163").Value = 468
$ws.Range("K163").Value = 1547
$ws.Range("L163").Value = 292.5
$ws.Range("M163").Value = 281
$ws.Range("N163").Value = 135
$ws.Range("O163").Value = 363
$ws.Range("P163").Value = 179

$ws.Range("A164").Value = "25-08-2021"
$ws.Range("B164").Value = 314.3
$ws.Range("C164").Value = 385
$ws.Range("D164").Value = 199
$ws.Range("E164").Value = 303
$ws.Range("F164").Value = 90
$ws.Range("G164").Value = 119.4
$ws.Range("H164").Value = 13.8
$ws.Range("I164").Value = 153.6
$ws.Range("J164").Value = 460.4
$ws.Range("K164").Value = 1537.7
$ws.Range("L164").Value = 287.8
$ws.Range("M164").Value = 273
$ws.Range("N164").Value = 131
$ws.Range("O164").Value = 358
$ws.Range("P164").Value = 174

$ws.Range("A165").Value = "26-08-2021"
$ws.Range("B165").Value = 316.4
$ws.Range("C165").Value = 387
$ws.Range("D165").Value = 202
$ws.Range("E165").Value = 304
$ws.Range("F165").Value = 89.8
$ws.Range("G165").Value = 119.4
$ws.Range("H165").Value = 12.2
$ws.Range("I165").Value = 156.1
$ws.Range("J165").Value = 460.4
$ws.Range("K165").Value = 1541.8
$ws.Range("L165").Value = 289.4
$ws.Range("M165").Value = 273
$ws.Range("N165").Value = 134
$ws.Range("O165").Value = 359
$ws.Range("P165").Value = 178

$ws.Range("A166").Value = "27-08-2021"
$ws.Range("B166").Value = 317.1
$ws.Range("C166").Value = 386
$ws.Range("D166").Value = 203
$ws.Range("E166").Value = 306
$ws.Range("F166").Value = 89.90000000000001
$ws.Range("G166").Value = 119.1
$ws.Range("H166").Value = 13.9
$ws.Range("I166").Value = 158.8
$ws.Range("J166").Value = 461.3
$ws.Range("K166").Value = 1538.1
$ws.Range("L166").Value = 287.1
$ws.Range("M166").Value = 271
$ws.Range("N166").Value = 136
$ws.Range("O166").Value = 357
$ws.Range("P166").Value = 178

$ws.Range("A167").Value = "30-08-2021"
$ws.Range("B167").Value = 316.4
$ws.Range("C167").Value = 382
$ws.Range("D167").Value = 203
$ws.Range("E167").Value = 308
$ws.Range("F167").Value = 90.09999999999999
$ws.Range("G167").Value = 118.8
$ws.Range("H167").Value = 16.9
$ws.Range("I167").Value = 160.3
$ws.Range("J167").Value = 462.8
$ws.Range("K167").Value = 1508.6
$ws.Range("L167").Value = 282.6
$ws.Range("M167").Value = 269
$ws.Range("N167").Value = 136
$ws.Range("O167").Value = 352
$ws.Range("P167").Value = 177

$ws.Range("A168").Value = "31-08-2021"
$ws.Range("B168").Value = 313
$ws.Range("C168").Value = 381
$ws.Range("D168").Value = 199
$ws.Range("E168").Value = 303
$ws.Range("F168").Value = 89.8
$ws.Range("G168").Value = 117
$ws.Range("H168").Value = 19.1
$ws.Range("I168").Value = 157.2
$ws.Range("J168").Value = 455.9
$ws.Range("K168").Value = 1512.6
$ws.Range("L168").Value = 282.5
$ws.Range("M168").Value = 272
$ws.Range("N168").Value = 136
$ws.Range("O168").Value = 352
$ws.Range("P168").Value = 175

$ws.Range("A169").Value = "01-09-2021"
$ws.Range("B169").Value = 312.9
$ws.Range("C169").Value = 381
$ws.Range("D169").Value = 200
$ws.Range("E169").Value = 303
$ws.Range("F169").Value = 89.40000000000001
$ws.Range("G169").Value = 117.7
$ws.Range("H169").Value = 14.6
$ws.Range("I169").Value = 157.6
$ws.Range("J169").Value = 458.1
$ws.Range("K169").Value = 1493
$ws.Range("L169").Value = 282.8
$ws.Range("M169").Value = 271
$ws.Range("N169").Value = 136
$ws.Range("O169").Value = 355
$ws.Range("P169").Value = 175

$ws.Range("A170").Value = "02-09-2021"
$ws.Range("B170").Value = 312.8
$ws.Range("C170").Value = 381
$ws.Range("D170").Value = 199
$ws.Range("E170").Value = 304
$ws.Range("F170").Value = 88.59999999999999
$ws.Range("G170").Value = 117
$ws.Range("H170").Value = 12.5
$ws.Range("I170").Value = 156.4
$ws.Range("J170").Value = 461.6
$ws.Range("K170").Value = 1481.6
$ws.Range("L170").Value = 283.2
$ws.Range("M170").Value = 271
$ws.Range("N170").Value = 137
$ws.Range("O170").Value = 354
$ws.Range("P170").Value = 178

$ws.Range("A171").Value = "03-09-2021"
$ws.Range("B171").Value = 310.4
$ws.Range("C171").Value = 379
$ws.Range("D171").Value = 196
$ws.Range("E171").Value = 302
$ws.Range("F171").Value = 88.2
$ws.Range("G171").Value = 115.7
$ws.Range("H171").Value = 15.5
$ws.Range("I171").Value = 154.3
$ws.Range("J171").Value = 460.6
$ws.Range("K171").Value = 1486.7
$ws.Range("L171").Value = 282.8
$ws.Range("M171").Value = 270
$ws.Range("N171").Value = 135
$ws.Range("O171").Value = 350
$ws.Range("P171").Value = 175

$dateRange.Style = "Normal"
